# Update of cost tables
# Applies the edits to the "costs (U.S. Cust.)" worksheet:
#  - Excavate/fill alluvial material unit cost (D8) becomes a constant, and
#    its total-cost formula (G8) is simplified to use MAX only.
#  - Angular boulder placement unit cost (D16) becomes a constant.
#  - Roads: Develop existing quantity (D42) doubles.
#  - A new "FEES AND LICENSING" section replaces "ENGINEERING FEES":
#      * "Markups (overhead, profit, insurance) and Engineering fees" (was
#        "From total costs") with an updated literature source and a new
#        16.5% rate.
#      * A new "Permitting" line at 35% of construction costs.
#  - The "TOTAL COSTS" row moves down one row and its formula now also
#    includes the new Permitting line.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("costs (U.S. Cust.)")

# --- Terraforming: Excavate/fill alluvial material -----------------------
$ws.Range("D8").Value = 23
$ws.Range("G8").Formula = "=D8*MAX(terraforming_volumes!C5,terraforming_volumes!C6)"

# --- Plant-stabilizing bioengineering: Angular boulder placement ---------
$ws.Range("D16").Value = 150

# --- Civil engineering & other: Roads - Develop existing ------------------
$ws.Range("D42").Value = 100

# --- Insert a new row for "Permitting" under the fees section ------------
# Row 59 is currently blank (between "From total costs" on row 58 and
# "TOTAL COSTS" on row 60); inserting here pushes "TOTAL COSTS" to row 61
# and copies formatting down from row 58, which is what we want for the new
# "Permitting" line.
$ws.Rows.Item(59).Insert()

# --- Rename section header: ENGINEERING FEES -> FEES AND LICENSING -------
$ws.Range("B57").Value = "FEES AND LICENSING"

# --- Row 58: "From total costs" -> "Markups ... and Engineering fees" ----
$ws.Range("B58").Value = "Markups (overhead, profit, insurance) and Engineering fees"
$ws.Range("D58").NumberFormat = "0.000"
$ws.Range("D58").Value = 0.165
$ws.Range("I58").Value = "LCH (2012)`nCramer (2012)`nJohnson (2019)"
$ws.Rows.Item(58).RowHeight = 48.95

# --- New row 59: "Permitting" --------------------------------------------
$ws.Range("B59:C59").Merge()
$ws.Range("B59").Value = "Permitting"
$ws.Range("D59").Value = 0.35
$ws.Range("E59").Value = "[-]"
$ws.Range("F59").Value = 1
$ws.Range("G59").Formula = "=G50*D59"
$ws.Rows.Item(59).RowHeight = 16.5

# --- Row 61 (was 60): TOTAL COSTS formula now includes Permitting --------
$ws.Range("G61").Formula = "=G58+G56+G55+G50+G59"

# --- Update the top-level return-on-investment formula on row 2 ----------
$ws.Range("G2").Formula = "=G61"

# --- Restore the view: scrolled down a bit, with G58 selected ------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("G58").Select()
